$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last data row (row 30) - the debtor record that was removed upon upload
$ws.Rows("30").Delete()

# Update remaining rows (2-29) with the new client names, dates and amounts
$ws.Range("B2").Value = "ABEL POLLO"
$ws.Range("C2").Value = 46046
$ws.Range("D2").Value = 83000

$ws.Range("B3").Value = "ALISO"
$ws.Range("C3").Value = 46048
$ws.Range("D3").Value = 100000

$ws.Range("B4").Value = "CAMILIN"
$ws.Range("C4").Value = 46046
$ws.Range("D4").Value = 545000

$ws.Range("B5").Value = "CAMPO VERDE TOCANCIPA"
$ws.Range("C5").Value = 46044
$ws.Range("D5").Value = 285000

$ws.Range("B6").Value = "CARNES JOHANA"
$ws.Range("C6").Value = 46045
$ws.Range("D6").Value = 164000

$ws.Range("B7").Value = "CARNILANDIA"
$ws.Range("C7").Value = 46046
$ws.Range("D7").Value = 280000

$ws.Range("B8").Value = "CIMARRON DORADO"
$ws.Range("C8").Value = 46041
$ws.Range("D8").Value = 426500

$ws.Range("B9").Value = "CIMARRON DORADO"
$ws.Range("C9").Value = 46045
$ws.Range("D9").Value = 403800

$ws.Range("B10").Value = "COCINA CHINA"
$ws.Range("C10").Value = 46047
$ws.Range("D10").Value = 170000

$ws.Range("B11").Value = "CRISTIAN ACACIAS"
$ws.Range("C11").Value = 46042
$ws.Range("D11").Value = 1000000

$ws.Range("B12").Value = "DARWIN FUTBOL"
$ws.Range("C12").Value = 45921
$ws.Range("D12").Value = 200000

$ws.Range("B13").Value = "DAVIDCITO"
$ws.Range("C13").Value = 45947
$ws.Range("D13").Value = 100000

$ws.Range("B14").Value = "DOÑA SANDRA"
$ws.Range("C14").Value = 46039
$ws.Range("D14").Value = 100000

$ws.Range("B15").Value = "EL RUBY"
$ws.Range("C15").Value = 46045
$ws.Range("D15").Value = 188000

$ws.Range("B16").Value = "FRANCO"
$ws.Range("C16").Value = 45996
$ws.Range("D16").Value = 20000

$ws.Range("B17").Value = "JUAN DAVID"
$ws.Range("C17").Value = 46046
$ws.Range("D17").Value = 560000

$ws.Range("B18").Value = "JULIANA POLLO"
$ws.Range("C18").Value = 46042
$ws.Range("D18").Value = 264000

$ws.Range("B19").Value = "LA SELECTA"
$ws.Range("C19").Value = 45912
$ws.Range("D19").Value = 82000

$ws.Range("B20").Value = "MERKA FRUVER DEXI"
$ws.Range("C20").Value = 45995
$ws.Range("D20").Value = 339000

$ws.Range("B21").Value = "MERKA FRUVER DEXI"
$ws.Range("C21").Value = 45988
$ws.Range("D21").Value = 15400

$ws.Range("B22").Value = "NEVADA"
$ws.Range("C22").Value = 46031
$ws.Range("D22").Value = 21900

$ws.Range("B23").Value = "NEVADA"
$ws.Range("C23").Value = 46038
$ws.Range("D23").Value = 175800

$ws.Range("B24").Value = "PARAISO FUNZA"
$ws.Range("C24").Value = 46045
$ws.Range("D24").Value = 173000

$ws.Range("B25").Value = "PINILLA"
$ws.Range("C25").Value = 45931
$ws.Range("D25").Value = 82000

$ws.Range("B26").Value = "PINILLA SOACHA"
$ws.Range("C26").Value = 46046
$ws.Range("D26").Value = 148000

$ws.Range("B27").Value = "PLAZA JESSICA"
$ws.Range("C27").Value = 46047
$ws.Range("D27").Value = 1344000

$ws.Range("B28").Value = "PREMIUM"
$ws.Range("C28").Value = 46046
$ws.Range("D28").Value = 178000

$ws.Range("B29").Value = "PUNTA DE ANCA"
$ws.Range("C29").Value = 46045
$ws.Range("D29").Value = 965000
